# Generate Report for Handoff
# Regenerate the localization-status report: new source-doc GUID
# (442560d1-f5a0-4978-91ba-be8a756b9074 -> 3ae5b4a1-f782-4e07-859a-f90961d9d7f3),
# new handoff-xliff content hash (e6e838ebfac7a1d1aacd09a266578d719d9f1d66 ->
# d05692d9e3904b7f8bc0ed4efe5f27ca373d7998), and refreshed timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "442560d1-f5a0-4978-91ba-be8a756b9074"
$newGuid = "3ae5b4a1-f782-4e07-859a-f90961d9d7f3"
$oldHash = "e6e838ebfac7a1d1aacd09a266578d719d9f1d66"
$newHash = "d05692d9e3904b7f8bc0ed4efe5f27ca373d7998"

$oldGenerateDate  = "2016-09-06 01:03:57"
$newGenerateDate  = "2016-09-06 01:04:20"
$oldHandoffDate   = "2016-09-06 01:03:52"
$newHandoffDate   = "2016-09-06 01:04:16"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2831fe9f14922555652074506615ba4ae06e710d/e2e/$oldGuid.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# The workbook's original hyperlink font: Calibri 11, single underline,
# cornflower-blue (RGB 0x6495ED == OLE/BGR 15570276). Re-applied after every
# Hyperlinks.Add below so re-creating the link doesn't silently swap the
# cell onto Excel's theme-colored built-in "Hyperlink" style.
function Restyle-HyperlinkCell($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

# --- Overview sheet ---
# A2: source file name
$wsOverview.Range("A2").Value = "$newGuid.md"

# B2: path and name (also the hyperlink's display text)
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md") | Out-Null
Restyle-HyperlinkCell $wsOverview.Range("B2")

# G2: latest HO xliff generate date
$wsOverview.Range("G2").Value = $newGenerateDate

# --- zh-cn sheet ---
# A2: source file name (hyperlink display text)
$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null
Restyle-HyperlinkCell $wsZhCn.Range("A2")

# G2: latest handoff file name
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"

# H2: latest handoff datetime
$wsZhCn.Range("H2").Value = $newHandoffDate

# --- de-de sheet ---
# A2: source file name (hyperlink display text)
$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null
Restyle-HyperlinkCell $wsDeDe.Range("A2")

# G2: latest handoff file name
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"

# H2: latest handoff datetime (shares the same value as Overview!G2)
$wsDeDe.Range("H2").Value = $newGenerateDate
